$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.137151122093201
$ws.Range("B1").Value = 3.209794282913208
$ws.Range("C1").Value = 3.562406301498413
$ws.Range("D1").Value = 3.942654132843018
$ws.Range("E1").Value = 1.207491874694824
